$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / label text updates
$ws.Range("A1").Value = "Sensor 1"
$ws.Range("B3").Value = "Dato Min: 10"
$ws.Range("B4").Value = "Dato Actual: 120"
$ws.Range("B5").Value = "Dato Max: 90"

# Numeric data updates
$ws.Range("A3").Value = 50
$ws.Range("C3").Value = 500
$ws.Range("A4").Value = 30
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 57
$ws.Range("A7").Value = 57
$ws.Range("A8").Value = 57
$ws.Range("A9").Value = 2000

# New row 10 - mirror the numeric style used by A6:A9
$ws.Range("A10").Value = 120
$ws.Range("A10").Font.Color = $ws.Range("A9").Font.Color
